$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: relocate the two rich-text note cells first (preserve their formatted runs) ---
# old D3 (DAILY note) -> new E4 ; old D4 (MONTHLY note) -> new E5
# (Cut automatically empties the source cell, so D3/D4 are clear afterwards.)
$ws.Range("D3").Cut($ws.Range("E4"))
$ws.Range("D4").Cut($ws.Range("E5"))

# --- Step 2: wipe everything else so the table can be rebuilt cleanly, without touching E4:E5 ---
$ws.Range("A1:D17").Clear()
$ws.Range("E1:E3").Clear()
$ws.Range("E6:E17").Clear()
$ws.Range("F1:F17").Clear()

# --- Step 3: column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.5546875
$ws.Columns.Item(2).ColumnWidth = 21.44140625
$ws.Columns.Item(3).ColumnWidth = 15.5546875
$ws.Columns.Item(4).ColumnWidth = 9.5546875

# --- Step 4: header row ---
$ws.Range("A1").Value = "Date "
$ws.Range("B1").Value = "Question/Analysis"
$ws.Range("C1").Value = "Scenario Name"
$ws.Range("D1").Value = "Models"
$ws.Range("E1").Value = "Inputs"
$ws.Range("F1").Value = "Notes"

# --- Step 5: row 2 (up2018 baseline) ---
$ws.Range("A2").Value = 43617
$ws.Range("A2").NumberFormat = "mmm-yy"
$ws.Range("C2").Value = "up2018"
$ws.Range("D2").Value = "SWBM"
$ws.Range("D2").NumberFormat = "mmm-yy"
$ws.Range("E2").Value = "Averaged cal-FJ precip, NA = 0; original Eto through Sep 2011, then spatial CIMIS."
$ws.Range("F2").Value = "GW and SW irrigation lower than expected in 2011-2018. "

# --- Step 6: row 3 (question box top) ---
$ws.Range("A3").Value = 43692
$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = "Question: Will daily vs monthly ET values make a difference in the water budget? "

# --- Step 7: row 4 (daily scenario) ---
$ws.Range("A4").Value = 43696
$ws.Range("A4").NumberFormat = "mm-dd-yy"
$ws.Range("C4").Value = "up2018_a"
$ws.Range("D4").Value = "SWBM"
$ws.Range("D4").NumberFormat = "mm-dd-yy"
# E4 already holds the DAILY rich-text note from the move above

# --- Step 8: row 5 (monthly scenario) ---
$ws.Range("A5").Value = 43696
$ws.Range("A5").NumberFormat = "mm-dd-yy"
$ws.Range("C5").Value = "up2018_b"
$ws.Range("D5").Value = "SWBM"
$ws.Range("D5").NumberFormat = "mm-dd-yy"
# E5 already holds the MONTHLY rich-text note from the move above

# --- Step 9: row 6 (question box bottom / result) ---
$ws.Range("A6").Value = 43696
$ws.Range("A6").NumberFormat = "mm-dd-yy"
$ws.Range("B6").Value = "Result: No visual difference between daily and monthly ET. Overall budgets start to diverge on the 3rd or 4th sigfig. Might as well use the monthly values. "

# --- Step 10: row 7 (analysis note) ---
$ws.Range("A7").Value = 43696
$ws.Range("A7").NumberFormat = "mm-dd-yy"
$ws.Range("B7").Value = "Analysis: How different will the 9 ECI273 scenarios make this water budget?"

# --- Step 11: row 8 (historical scenario) ---
$ws.Range("A8").NumberFormat = "mm-dd-yy"
$ws.Range("C8").Value = "hist"
$ws.Range("D8").Value = "SWBM"
$ws.Range("E8").Value = "Historical precip, gaps filled with ranked regression, created 2019.08.19 (leapdays now included!). Monthly ET (up2018_b). "

# --- Step 12: rows 9-17 (the 9 ECI273 precip-variant scenarios) ---
$ws.Range("C9").Value = "pvar_a10"
$ws.Range("D9").Value = "SWBM"
$ws.Range("C10").Value = "pvar_a5"
$ws.Range("C11").Value = "pvar_a3"
$ws.Range("C12").Value = "pvar_b90"
$ws.Range("C13").Value = "pvar_b80"
$ws.Range("C14").Value = "pvar_b70"
$ws.Range("C15").Value = "pvar_c10"
$ws.Range("C16").Value = "pvar_c20"
$ws.Range("C17").Value = "pvar_c30"

# --- Step 13: the medium box border surrounding the Question/Analysis block A3:F6 ---
$ws.Range("A3:F6").BorderAround($null, -4138)

# --- Step 14: row heights on rows 2 and 6 (thick-bottom rows grow slightly) ---
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15

# --- Step 15: selection / view state ---
$ws.Range("E9").Select()

Write-Host "done"
